# Generate Report for Handback
# Populate the "e7d3664c-2171-4c03-a1fc-289786d4cc8d" handback row (row 5) with the
# (failed) handback-validation results on both locale sheets (zh-cn, de-de), and
# widen the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9f63959e98476c9f42c81f737a398ba18969ec5/e2e/e7d3664c-2171-4c03-a1fc-289786d4cc8d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f302e1af6a422ce2d2f00fbbea0f5486d5a2d5e/e2e/e7d3664c-2171-4c03-a1fc-289786d4cc8d.md."
$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f302e1af6a422ce2d2f00fbbea0f5486d5a2d5e/e2e/e7d3664c-2171-4c03-a1fc-289786d4cc8d.md"
$mdDisplay = "e7d3664c-2171-4c03-a1fc-289786d4cc8d.md"

function Update-LocaleSheet($ws, $xlfName, $handbackDateTime) {
    # Widen the "Error Detail" column (16 / P) so the new message is legible.
    $refWidth = $ws.Columns.Item(1).ColumnWidth
    $ws.Columns.Item(16).ColumnWidth = $refWidth

    # I5 - "Latest Target File": now resolved, points at the source .md file.
    $ws.Range("I5").Value = $mdDisplay
    $ws.Hyperlinks.Add($ws.Range("I5"), $latestMdUrl, $null, $null, $mdDisplay) | Out-Null
    $ws.Range("I5").Font.Underline = 2
    $ws.Range("I5").Font.Color = 15570276

    # J5 - "Latest Handback File": the generated xlf for this locale.
    $ws.Range("J5").Value = $xlfName

    # K5 - "Latest Handback DateTime": when the handback xlf was generated.
    $ws.Range("K5").Value = $handbackDateTime

    # P5 - "Error Detail": report that the handback wasn't built off the latest commit.
    $ws.Range("P5").Value = $errorMessage
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LocaleSheet $wsZhCn "e7d3664c-2171-4c03-a1fc-289786d4cc8d.0e50bb54dbe8a40e425dad31745244f29197b78a.zh-cn.xlf" "2016-10-17 16:03:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LocaleSheet $wsDeDe "e7d3664c-2171-4c03-a1fc-289786d4cc8d.0e50bb54dbe8a40e425dad31745244f29197b78a.de-de.xlf" "2016-10-17 16:03:40"
